$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.45"
$ws.Range("E2").Value = "'1.35%"
$ws.Range("D3").Value = "'41.29"
$ws.Range("E3").Value = "'0.31%"
$ws.Range("D4").Value = "'5.134"
$ws.Range("E4").Value = "'0.60%"
$ws.Range("D5").Value = "'0.07707"
$ws.Range("E5").Value = "'1.37%"
$ws.Range("D6").Value = "'1.626"
$ws.Range("E6").Value = "'0.21%"
$ws.Range("D7").Value = "'0.9247"
$ws.Range("E7").Value = "'2.02%"
$ws.Range("D9").Value = "'0.1203"
$ws.Range("E9").Value = "'18.02%"
$ws.Range("D10").Value = "'0.1838"
$ws.Range("E10").Value = "'4.61%"
$ws.Range("E11").Value = "'0.40%"
$ws.Range("E12").Value = "'-0.35%"
$ws.Range("E13").Value = "'-0.51%"
$ws.Range("D14").Value = "'0.001251"
$ws.Range("E14").Value = "'1.43%"
$ws.Range("D15").Value = "'0.005767"
$ws.Range("E15").Value = "'-1.81%"
$ws.Range("D16").Value = "'3.353"
$ws.Range("E16").Value = "'0.07%"
$ws.Range("D17").Value = "'4.312"
$ws.Range("E17").Value = "'1.16%"
$ws.Range("E18").Value = "'1.89%"
$ws.Range("D19").Value = "'6.937"
$ws.Range("E19").Value = "'6.00%"
$ws.Range("D20").Value = "'0.1389"
$ws.Range("E20").Value = "'2.45%"
$ws.Range("D22").Value = "'0.04057"
$ws.Range("E22").Value = "'-3.01%"
$ws.Range("E23").Value = "'2.91%"
$ws.Range("D24").Value = "'0.004112"
$ws.Range("E24").Value = "'1.23%"
$ws.Range("E25").Value = "'-2.43%"
$ws.Range("E26").Value = "'24.69%"
$ws.Range("E38").Value = "'3.62%"
$ws.Range("D39").Value = "'0.05279"
$ws.Range("E39").Value = "'2.48%"
$ws.Range("D40").Value = "'0.007839"
$ws.Range("E40").Value = "'0.75%"
$ws.Range("D41").Value = "'0.1316"
$ws.Range("E41").Value = "'1.52%"
$ws.Range("D42").Value = "'0.006794"
$ws.Range("E42").Value = "'-2.79%"
$ws.Range("D43").Value = "'0.001844"
$ws.Range("E43").Value = "'-3.93%"
$ws.Range("D44").Value = "'0.008181"
$ws.Range("E44").Value = "'-3.23%"
$ws.Range("D45").Value = "'0.3096"
$ws.Range("E45").Value = "'-7.10%"
$ws.Range("D46").Value = "'0.00006732"
$ws.Range("E46").Value = "'5.93%"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("D48").Value = "'0.2055"
$ws.Range("E48").Value = "'2,176.81%"
$ws.Range("E49").Value = "'-6.93%"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("E51").Value = "'-0.13%"
